# Update LLBv3 Pin Mapping sheet to match current state of board.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "LLBV3 Header" -> "LLBV3 Header / Function"
$ws.Range("E1").Value = "LLBV3 Header / Function"

# Fill in the "Header / Function" column (E) for pins that previously had
# a Pin Name (D) but no header/function note.
$ws.Range("E2").Value  = "MCP 2515 interrupt on received frames"
$ws.Range("E6").Value  = "X3, for power on board"
$ws.Range("E7").Value  = "X3, for power on board"
$ws.Range("E8").Value  = "E-stop jumper, also X3"
$ws.Range("E18").Value = "X3, for power on board"
$ws.Range("E19").Value = "wheel hall switch header"
$ws.Range("E21").Value = "all SPI devices, SPI header"
$ws.Range("E22").Value = "all SPI devices, SPI header"
$ws.Range("E23").Value = "all SPI devices, SPI header"
$ws.Range("E20").Value = "SPI header (this pin tells the mega to be a slave)"
$ws.Range("E36").Value = "MCP2515 slave selection"
$ws.Range("E37").Value = "DAC slave selection"
$ws.Range("E53").Value = "on-board buzzer"
$ws.Range("E55").Value = "X3, no purpose assgined"
$ws.Range("E57").Value = "X3, for power on board"

# Match author's last-saved view state (selection ended up on E58).
[void]$ws.Range("E58").Select()
